$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# --- Update the three Cypher query cells (B2:B4) ---

# Row 2: CasesTab detail query -- append an ORDER BY / LIMIT clause
$ws.Range("B2").Value = "MATCH (ss:study_subject)" + [char]10 + `
"MATCH (ss)<-[:sample_of_study_subject]-(sp)<-[:file_of_sample]-(f)-[:file_of_laboratory_procedure]->(lp)" + [char]10 + `
"WITH ss, collect(DISTINCT sp.sample_id) AS samples, collect(DISTINCT lp.laboratory_procedure_id) AS lab_procedures, collect(DISTINCT f) AS files" + [char]10 + `
"MATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)" + [char]10 + `
"MATCH (ss)<-[:sf_of_study_subject]-(sf)" + [char]10 + `
"MATCH (ss)<-[:diagnosis_of_study_subject]-(d)" + [char]10 + `
"MATCH (d)<-[:tp_of_diagnosis]-(tp)" + [char]10 + `
"MATCH (ss)<-[:demographic_of_study_subject]-(demo)" + [char]10 + `
" WHERE   tp.endocrine_therapy_type IN [`"None`"]" + [char]10 + `
"return ss.study_subject_id as ``Case ID``," + [char]10 + `
"       p.program_acronym as ``Program Code``," + [char]10 + `
"        p.program_id as Program_ID," + [char]10 + `
"       s.study_acronym as ``Arm``," + [char]10 + `
"       ss.disease_subtype as ``Diagnosis``," + [char]10 + `
"       sf.grouped_recurrence_score AS ``Recurrence Score``," + [char]10 + `
"       d.tumor_size_group AS ``tumor_size``," + [char]10 + `
"       d.er_status AS ``ER Status``," + [char]10 + `
"       d.pr_status AS ``PR Status``," + [char]10 + `
"       coalesce(CASE demo.age_at_index % 1 WHEN 0 THEN apoc.convert.toInteger(demo.age_at_index) ELSE demo.age_at_index END, '') AS ``Age (years)``," + [char]10 + `
"demo.survival_time AS ``Survival (days)``" + [char]10 + `
" order By ss.study_subject_id ASC LIMIT 100 "

# Row 3: SamplesTab query -- append an ORDER BY / LIMIT clause
$ws.Range("B3").Value = "MATCH (ss:study_subject)" + [char]10 + `
"WITH COLLECT(ss.study_subject_id) AS all_subjects" + [char]10 + `
"MATCH (samp:sample)" + [char]10 + `
"MATCH (samp)-[:sample_of_study_subject]->(ss)" + [char]10 + `
"MATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)" + [char]10 + `
"MATCH (samp)<-[:file_of_sample]-(f)-[:file_of_laboratory_procedure]->(lp)" + [char]10 + `
"MATCH (ss)<-[:diagnosis_of_study_subject]-(d)" + [char]10 + `
"MATCH (d)<-[:tp_of_diagnosis]-(tp)" + [char]10 + `
" WHERE   tp.endocrine_therapy_type IN [`"None`"] " + [char]10 + `
"WITH" + [char]10 + `
"    distinct lp," + [char]10 + `
"    toInteger(split(ss.study_subject_id,'-')[2]) AS subject_id_num," + [char]10 + `
"    collect(distinct f.file_id) AS files," + [char]10 + `
"    samp, ss, s, p, all_subjects" + [char]10 + `
"RETURN" + [char]10 + `
" samp.sample_id AS ``Sample ID``," + [char]10 + `
"            ss.study_subject_id AS ``Case ID``," + [char]10 + `
"            p.program_acronym AS ``Program Code``," + [char]10 + `
"            s.study_acronym AS ``Arm``," + [char]10 + `
"            ss.disease_subtype AS ``Diagnosis``," + [char]10 + `
"            samp.tissue_type AS ``Tissue Type``," + [char]10 + `
"            samp.composition AS ``Tissue Composition``," + [char]10 + `
"            samp.sample_anatomic_site AS ``Sample Anatomic Site``," + [char]10 + `
"            samp.method_of_sample_procurement AS ``Sample Procurement Method``" + [char]10 + `
" order By samp.sample_id ASC LIMIT 100"

# Row 4: FilesTab query -- change "order by" to "order By ... LIMIT 100"
$ws.Range("B4").Value = "MATCH (f:file)-->(parent)" + [char]10 + `
"MATCH (f)-[:file_of_sample]->(samp)" + [char]10 + `
"MATCH (samp)-[:sample_of_study_subject]->(ss)" + [char]10 + `
"MATCH (ss)-[:study_subject_of_study]->(s)" + [char]10 + `
"MATCH (s)-[:study_of_program]->(p)" + [char]10 + `
"MATCH (d)-[:diagnosis_of_study_subject]->(ss)" + [char]10 + `
"MATCH (tp)-[:tp_of_diagnosis]->(d)" + [char]10 + `
" WHERE   tp.endocrine_therapy_type IN [`"None`"] " + [char]10 + `
"WITH" + [char]10 + `
"        f, parent,p, ss, d,tp, s, samp," + [char]10 + `
"        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units," + [char]10 + `
"        toInteger(floor(log(f.file_size)/log(1024))) as i," + [char]10 + `
"        2 as precision" + [char]10 + `
"WITH" + [char]10 + `
"        f, parent,p, ss, d,tp, s, samp," + [char]10 + `
"        f.file_size /(1024^i) AS value," + [char]10 + `
"        10^precision AS factor," + [char]10 + `
"        units[i] as unit" + [char]10 + `
"WITH" + [char]10 + `
"        f, parent,p, ss, d,tp, s, samp, unit," + [char]10 + `
"        round(factor * value)/factor AS size" + [char]10 + `
"RETURN Distinct" + [char]10 + `
"    f.file_name AS ``File Name``," + [char]10 + `
"    head(labels(samp)) AS ``Association``," + [char]10 + `
"    f.file_description AS ``Description``," + [char]10 + `
"    f.file_format AS ``File Format``," + [char]10 + `
"     CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size," + [char]10 + `
"    p.program_acronym AS ``Program Code``," + [char]10 + `
"    s.study_acronym AS ``Arm``," + [char]10 + `
"    ss.study_subject_id AS ``Case ID``," + [char]10 + `
"    samp.sample_id AS ``Sample ID``" + [char]10 + `
"    order By f.file_name ASC LIMIT 100"

# --- Row heights grow because the cell text got longer (autofit wrapped text) ---
$ws.Rows.Item(2).RowHeight = 345.6
$ws.Rows.Item(3).RowHeight = 360
$ws.Rows.Item(4).RowHeight = 409.6

# --- Update sheet view / selection state ---
$ws.Range("B4").Select()
$excel.ActiveWindow.ScrollRow = 3
